$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("B11").Value = "Schneid, was du Lust hast."
$ws.Range("C11").Value = "Ich habe meiner Coiffeuse erlaubt, das zu schneiden, was sie Lust hat und sie findet, würde noch passen. Sie hat dann letztlich nicht wirklich variiert, sondern nur was kleines verändert. Meine Frau war abends erstaunt, als ich es ihr erzählt habe. So etwas hätte ich früher nie gemacht. Und ja, stimmt maximal - das hätte ich früher definitiv nicht gemacht."
$ws.Range("D11").Value = "Mut"
$ws.Range("E11").Value = "https://www.biblond.com/wp-content/uploads/2014/11/biblond_web_magazine_coiffure_pro_tendance-coiffure_style-hipster_cyril-hohl-3-800x600.jpg"

$ws.Range("E12").Select()
